# faturamento_diario.xlsx update:
#  - correct 3 existing "05/2025" daily totals (days 12, 19, 20)
#  - insert 3 new daily rows for "05/2025" (days 21, 22, 23) right after day 20,
#    pushing every following row down by three

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix existing values for days 12, 19 and 20 (rows 8, 13, 14) ---
$ws.Cells.Item(8, 2).Value = 13948.05
$ws.Cells.Item(13, 2).Value = 15817.9
$ws.Cells.Item(14, 2).Value = 20392.32

# --- insert 3 fresh rows right before the old row 15 (day 1 / Abril) ---
$ws.Rows.Item(15).Resize(3, 1).Insert()

# --- populate the 3 new rows: days 21, 22, 23 of 05/2025 ---
$ws.Cells.Item(15, 1).Value = 21
$ws.Cells.Item(15, 2).Value = 7699.42
$ws.Cells.Item(15, 3).Value = 5
$ws.Cells.Item(15, 4).Value = 2025
$ws.Cells.Item(15, 5).Value = "05/2025"

$ws.Cells.Item(16, 1).Value = 22
$ws.Cells.Item(16, 2).Value = 27031.65
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 2025
$ws.Cells.Item(16, 5).Value = "05/2025"

$ws.Cells.Item(17, 1).Value = 23
$ws.Cells.Item(17, 2).Value = 6415.3
$ws.Cells.Item(17, 3).Value = 5
$ws.Cells.Item(17, 4).Value = 2025
$ws.Cells.Item(17, 5).Value = "05/2025"
